$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "Cases" query (row 2 / B2) loses its trailing "Cohort" return column
# (keeps everything through "Response to Treatment").
$casesQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`nMATCH (c)<--(diag:diagnosis)`nOPTIONAL MATCH (samp:sample)-->(c)`nOPTIONAL MATCH (co:cohort)<-[*]-(c)`nWITH DISTINCT c, s, demo, diag, co`nWHERE diag.stage_of_disease IN ['IIIa']`nRETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n        coalesce(demo.breed, '') AS Breed ,`n        coalesce(diag.disease_term, '') AS Diagnosis ,`n        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n        coalesce(demo.sex, '') AS Sex ,`n        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n        coalesce(demo.weight, '') AS ``Weight (kg)``,`n        coalesce(diag.best_response, '') AS ``Response to Treatment``"

# The "Samples" query (row 3 / B3) text, unchanged in content.
$sampleQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) `nWHERE diag.stage_of_disease IN ['IIIa']`nWITH DISTINCT samp AS samp, c, demo, diag`nRETURN  coalesce(samp.sample_id, '') AS ``Sample ID``, `n        coalesce(c.case_id, '') AS ``Case ID``, `n        coalesce(demo.breed,'') AS Breed , `n        coalesce(diag.disease_term,'') AS Diagnosis , `n        coalesce(samp.sample_site, '') AS ``Sample Site``,`n        coalesce(samp.summarized_sample_type, '') AS ``Sample Type``,`n        coalesce(samp.specific_sample_pathology, '') AS ``Pathology/Morphology``,`n        coalesce(samp.tumor_grade, '') AS ``Tumor Grade``,`n        coalesce(samp.sample_chronology, '') AS ``Sample Chronology``,`n        coalesce(samp.percentage_tumor, '') AS ``Percentage Tumor``,`n        coalesce(samp.necropsy_sample, '') AS ``Necropsy Sample``,`n        coalesce(samp.sample_preservation, '') AS ``Sample Preservation``"

# The "Files" query (row 4 / B4) text, unchanged in content.
$filesQuery = "`nMATCH (f:file)-->(parent)`nWITH DISTINCT f, parent`nMATCH (f)-[*]->(c:case)<--(demo:demographic)`n MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`nWHERE diag.stage_of_disease IN ['IIIa']`nWITH DISTINCT f, parent, c, demo, diag, s`nRETURN coalesce(f.file_name, '') AS ``File Name``, `n        coalesce(labels(parent)[0], '') AS ``Association``,`n        coalesce(f.file_description, '') AS ``Description``,`n        coalesce(f.file_format, '') AS ``Format``,`n        coalesce(f.file_size, '') AS ``Size``,`n        coalesce(c.case_id, '') AS ``Case ID``, `n        coalesce(diag.disease_term,'') AS Diagnosis , `n        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

$ws.Range("B2").Value = $casesQuery
$ws.Range("B3").Value = $sampleQuery
$ws.Range("B4").Value = $filesQuery

$ws.Range("B2").Select()
